# Auto-generated edit script applying the cryptos.xlsx (crypto price list) diff.
# All D (Price) and E (Volume) column values are plain text in this workbook
# (e.g. "51.795.22", "  +0.16%  "), so we force text NumberFormat before assigning
# to stop Excel auto-coercing number-looking strings into floating point numbers,
# then restore the default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '51.795.22'
Set-TextValue 'E2' '  +0.16%  '
Set-TextValue 'D3' '2.968.82'
Set-TextValue 'E3' '  +1.06%  '
Set-TextValue 'D4' '0.998'
Set-TextValue 'E4' '  -0.23%  '
Set-TextValue 'D5' '379.98'
Set-TextValue 'E5' '  +7.92%  '
Set-TextValue 'D6' '104.94'
Set-TextValue 'E6' '  -1.37%  '
Set-TextValue 'D7' '0.546'
Set-TextValue 'E7' '  -0.20%  '
Set-TextValue 'E8' '  -0.02%  '
Set-TextValue 'D9' '0.600'
Set-TextValue 'E9' '  +0.90%  '
Set-TextValue 'D10' '37.41'
Set-TextValue 'E10' '  +0.10%  '
Set-TextValue 'E11' '  -0.04%  '
Set-TextValue 'D12' '0.0844'
Set-TextValue 'E12' '  -0.06%  '
Set-TextValue 'D13' '18.62'
Set-TextValue 'E13' '  -0.69%  '
Set-TextValue 'D14' '3.425.88'
Set-TextValue 'E14' '  +0.68%  '
Set-TextValue 'D15' '7.48'
Set-TextValue 'E15' '  +1.01%  '
Set-TextValue 'D16' '2.939.41'
Set-TextValue 'E16' '  +0.09%  '
Set-TextValue 'D17' '0.960'
Set-TextValue 'E17' '  -1.71%  '
Set-TextValue 'D18' '51.698.91'
Set-TextValue 'E18' '  +0.27%  '
Set-TextValue 'D19' '3.48'
Set-TextValue 'E19' '  +5.23%  '
Set-TextValue 'D20' '7.43'
Set-TextValue 'E20' '  +2.29%  '
Set-TextValue 'D21' '13.17'
Set-TextValue 'E21' '  -0.13%  '
Set-TextValue 'D22' '0.0₃0959'
Set-TextValue 'E22' '  +0.63%  '
Set-TextValue 'D23' '68.73'
Set-TextValue 'E23' '  +0.03%  '
Set-TextValue 'D24' '263.61'
Set-TextValue 'E24' '  -0.19%  '
Set-TextValue 'D25' '2.80'
Set-TextValue 'E25' '  +5.00%  '
Set-TextValue 'D26' '7.46'
Set-TextValue 'E26' '  +20.06%  '
Set-TextValue 'B27' 'Kaspa'
Set-TextValue 'C27' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D27' '0.170'
Set-TextValue 'E27' '  -2.45%  '
Set-TextValue 'B28' 'LEO'
Set-TextValue 'C28' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D28' '4.16'
Set-TextValue 'E28' '  -3.99%  '
Set-TextValue 'E29' '  +0.02%  '
Set-TextValue 'D30' '7.43'
Set-TextValue 'E30' '  +3.29%  '
Set-TextValue 'D31' '26.04'
Set-TextValue 'E31' '  -1.38%  '
Set-TextValue 'E32' '  -4.00%  '
Set-TextValue 'D33' '9.93'
Set-TextValue 'E33' '  -0.73%  '
Set-TextValue 'D34' '52.76'
Set-TextValue 'E34' '  +3.95%  '
Set-TextValue 'D35' '34.44'
Set-TextValue 'E35' '  -2.49%  '
Set-TextValue 'E36' '  -4.07%  '
Set-TextValue 'D37' '0.0438'
Set-TextValue 'E37' '  +3.15%  '
Set-TextValue 'E38' '  +0.32%  '
Set-TextValue 'D39' '3.06'
Set-TextValue 'E39' '  -4.81%  '
Set-TextValue 'B40' 'Celestia'
Set-TextValue 'C40' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D40' '17.41'
Set-TextValue 'E40' '  +1.19%  '
Set-TextValue 'B41' 'Stacks'
Set-TextValue 'C41' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D41' '2.67'
Set-TextValue 'E41' '  -5.25%  '
Set-TextValue 'D42' '1.85'
Set-TextValue 'E42' '  -2.39%  '
Set-TextValue 'E43' '  +1.03%  '
Set-TextValue 'D44' '124.03'
Set-TextValue 'E44' '  +2.33%  '
Set-TextValue 'D45' '22.04'
Set-TextValue 'E45' '  -2.35%  '
Set-TextValue 'D46' '0.282'
Set-TextValue 'E46' '  +18.98%  '
Set-TextValue 'E47' '  -3.57%  '
Set-TextValue 'D48' '2.032.95'
Set-TextValue 'E48' '  -3.07%  '
Set-TextValue 'E49' '  +0.22%  '
Set-TextValue 'D50' '3.23'
Set-TextValue 'E50' '  -0.06%  '
Set-TextValue 'D51' '0.0332'
Set-TextValue 'E51' '  +4.12%  '

Write-Host "Applied 99 cell updates"
